# Adds a new dialogue row (row 5) to the Даскалл sheet, matching the
# "Add files via upload" commit: a new entry (line 257) with English,
# Russian and "encoded" Russian text, plus the matching separator border
# that moves from the old last row (row 4) onto the new last row boundary.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 becomes an interior row: its old "last row" bottom border is
# replaced by a plain thin bottom rule (this also creates the new border
# and the two new cellXfs entries used by A4:E4). ---
$row4 = $ws.Range("A4:E4")
$row4.Borders.Item(9).LineStyle = 1
$row4.Borders.Item(9).Weight = 2
$row4.Borders.Item(9).ColorIndex = -4105

# --- New row 5 content ---
$ws.Range("B5").Value = 257
$ws.Range("C5").Value = ' Startling, isn\''t it?[K] Visitations of\nPokémon from a distant future…'
$ws.Range("D5").Value = ' Поразительно, правда?[K] Покемоны,\nприбывшие из далёкого будущего...'
$ws.Range("E5").Value = ' Ðïñàèéóåìûîï, ðñàâäà?[K] Ðïëåíïîú,\nðñéáúâšéå éè äàìæëïãï áôäôþåãï…'

$ws.Rows.Item(5).RowHeight = 31.8

# --- Selection moves to C5, matching the saved workbook state ---
$ws.Range("C5").Select()
